$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mis-spelled "Amout" column header in the database schema sheet.
$ws.Range("I8").Value = "Amount"
